# Apply data corrections to vTiger_testdata.xlsx test script data.
$wb = $excel.ActiveWorkbook

# --- Org sheet ---
$ws = $wb.Worksheets.Item("Org")
$ws.Range("D5").Value = "Energy"
$ws.Range("E5").Value = "Press"
$ws.Range("C5").Value = "TKP_"
$ws.Range("C8").Value = "TekP_"

# --- Vendor sheet ---
$ws = $wb.Worksheets.Item("Vendor")
$ws.Range("C2").Value = "Google"

# --- Lead sheet ---
$ws = $wb.Worksheets.Item("Lead")
$ws.Range("E4").Value = "Web Site"
$ws.Range("F4").Value = "Education"

# --- Product sheet ---
$ws = $wb.Worksheets.Item("Product")
$ws.Range("C2").Value = "Books_"
$ws.Range("C5").Value = "Pen_"

# --- back to Org sheet for phone number ---
$ws = $wb.Worksheets.Item("Org")
$ws.Range("D8").Value = "'9760456233"

# --- Selections on other sheets ---
$wb.Worksheets.Item("Contact").Range("C5").Select()
$wb.Worksheets.Item("Vendor").Range("C8").Select()
$wb.Worksheets.Item("Lead").Range("C9").Select()
$wb.Worksheets.Item("Product").Range("D10").Select()

# --- re-activate Org sheet (keeps it the tab-selected sheet) and set its selection last ---
$wb.Worksheets.Item("Org").Activate()
$wb.Worksheets.Item("Org").Range("B10").Select()
